$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column F (PO_Spend -> PO_Commit)
$ws.Range("F1").Value = "PO_Commit"

# Update D2:D5 values from 40 to 0
$ws.Range("D2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("D5").Value = 0

# Update selection to E5
$ws.Range("E5").Select()
